# Actualización automática 2025-08-07 08:30:08
#
# LOZANO MOLINA TITO reported a new sale of 2227.24 (RENOVA&DISEÑA S.A.,
# category "PIEDRA SINTERIZADA") for the current month, so we push the
# figure into the three summary sheets and let the dependent totals /
# "x de N" counters / percentages follow.

$wb = $excel.ActiveWorkbook

# Helper: Excel's ColumnWidth property is expressed in "characters" but is
# internally snapped to a pixel grid, which makes it shift slightly away
# from the integer you assign (it bakes in the standard ~5 pixel cell
# padding). Feeding it (target - 11/12) lands back exactly on the desired
# stored width for this workbook's default font.
function Set-ExactColumnWidth($sheet, $colIndex, $targetWidth) {
    $sheet.Columns($colIndex).ColumnWidth = $targetWidth - (11/12)
}

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": per-client breakdown by product group
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# LOZANO MOLINA TITO / RENOVA&DISEÑA S.A., column L = "PIEDRA SINTERIZADA"
$wsGrupo.Range("L19").Value = 2227.24

# Row 29 tallies how many of the 27 clients have a non-zero value in each
# column; one more client (row 19) now has a value in column L.
$wsGrupo.Range("L29").Value = "1 de 27"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": per-client breakdown by month
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# LOZANO MOLINA TITO / RENOVA&DISEÑA S.A., column F = "agosto"
$wsMensual.Range("F19").Value = 2227.24

# Row 29 is the column total; add the same amount that was just booked.
$wsMensual.Range("F29").Value = 2227.24

# Column F widened by one character after the longer number was entered.
Set-ExactColumnWidth $wsMensual 6 13

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": budget vs. actual sales by product group
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 15 = "PIEDRA SINTERIZADA": VENTA, POR CUMPLIR (= PRESUPUESTO - VENTA)
# and CUMPLIMIENTO (= VENTA / PRESUPUESTO) all update together.
# (.Value2 is used for read-back; in this host .Value only works reliably
# as a setter.)
$wsCumpl.Range("D15").Value = 2227.24
$wsCumpl.Range("E15").Value = $wsCumpl.Range("C15").Value2 - $wsCumpl.Range("D15").Value2
$wsCumpl.Range("F15").Value = $wsCumpl.Range("D15").Value2 / $wsCumpl.Range("C15").Value2

# Row 19 = TOTAL row, same relationships.
$wsCumpl.Range("D19").Value = 2227.24
$wsCumpl.Range("E19").Value = $wsCumpl.Range("C19").Value2 - $wsCumpl.Range("D19").Value2
$wsCumpl.Range("F19").Value = $wsCumpl.Range("D19").Value2 / $wsCumpl.Range("C19").Value2

# Columns D (VENTA), E (POR CUMPLIR) and F (CUMPLIMIENTO) widened to fit
# the newly entered figures.
Set-ExactColumnWidth $wsCumpl 4 13
Set-ExactColumnWidth $wsCumpl 5 24
Set-ExactColumnWidth $wsCumpl 6 22

Write-Host "Applied 2025-08-07 08:30:08 update: LOZANO MOLINA TITO / RENOVA&DISENA S.A. PIEDRA SINTERIZADA = 2227.24"
